# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig —
# rebrand the StructureDefinition spreadsheet from "Alvearie"/ibm.com to
# "LinuxForHealth"/linuxforhealth.org, bump the version, and refresh the date.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: Property / Value table -------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/child-organization-hierarchy-level-code"
# Version
$wsMeta.Range("B3").Value = "8.0.0"
# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: StructureDefinition element table ---------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension") no longer carries the ele-1/ext-1 constraint text in the
# Constraint(s) column (it now only lives on the Extension.extension row).
$wsElem.Range("AI2").Value = ""

# Row 5 ("Extension.url") Fixed Value column: same URL rebrand as above.
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/child-organization-hierarchy-level-code"
